# SpellTrap.xlsx: "OnMagicDamage now can known who is the maker"
#
# The trap effect formula stored in D9 (row for 54000006 / "爆炸陷阱") is
# updated so that OnMagicDamage is called with an extra leading argument
# (the "maker") in addition to the existing target-damage/type arguments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "m.OnMagicDamage(null,t.Damage,3);return true;"

# Keep the active-cell selection in sync with where the edit was made,
# mirroring the saved view state after the change.
[void]$ws.Range("D9").Select()
